# Trade #3 closed at 2026-02-17 12:26:35 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.98
$summary.Range("B4").Value = -0.02
$summary.Range("B5").Value = -0.13
$summary.Range("B6").Value = 3
$summary.Range("B8").Value = 2
$summary.Range("B9").Value = 33.33

# ---------------------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.98
$status.Range("D4").Value = 3
$status.Range("E4").Value = -0.02
$status.Range("F4").Value = -0.02
$status.Range("G4").Value = 33.33

# ---------------------------------------------------------------------------
# Helper: append Trade #3 row to a trade-log sheet (All Trades / MarketMaking)
# ---------------------------------------------------------------------------
function Add-TradeThreeRow($sheet) {
    $sheet.Cells.Item(4, 1).Value = 3
    # Force the date column to stay plain text (matches the other rows,
    # which store the date as a literal string, not a date serial) -
    # otherwise Excel's COM layer auto-parses "2026-02-17" into a date.
    $sheet.Cells.Item(4, 2).NumberFormat = "@"
    $sheet.Cells.Item(4, 2).Value = "2026-02-17"
    $sheet.Cells.Item(4, 2).ClearFormats()
    $sheet.Cells.Item(4, 3).Value = "12:26:28"
    $sheet.Cells.Item(4, 4).Value = "MarketMaking"
    $sheet.Cells.Item(4, 5).Value = "UP"
    $sheet.Cells.Item(4, 6).Value = 0.19
    $sheet.Cells.Item(4, 7).Value = 0.18
    $sheet.Cells.Item(4, 8).Value = "CLOSED"
    $sheet.Cells.Item(4, 9).Value = -5.2632
    $sheet.Cells.Item(4, 10).Value = -0.01
    $sheet.Cells.Item(4, 11).Value = 99.98
    $sheet.Cells.Item(4, 12).Value = 0
    $sheet.Cells.Item(4, 13).Value = 0
    $sheet.Cells.Item(4, 14).Value = 0.6
    $sheet.Cells.Item(4, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item(4, 16).Value = "early_exit"
    $sheet.Cells.Item(4, 17).Value = 0.13
}

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeThreeRow $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeThreeRow $marketMaking
